$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" sheet
# and the "全部类型" sheet, which carry duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F7").Value = 135
    $ws.Range("F8").Value = 57
    $ws.Range("F9").Value = 368
}
